$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "Set up the node js server file to pull data from the backend and serve it to the frontend. Set up a communal email and MongoDB account for the team to use. Provided instructions for team to install npm backend modules (these are in instructions.txt in the gitbucket repo). Working with Cedric to pull data from the sign up form to the backend."

$ws.Range("D7").Select()
